$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.027.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.667.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5102'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06381'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07436'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.675.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.507'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5800'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008480'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.849.19'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.922'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.186'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.612'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1216'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.88%  '
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06617'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.13%  '
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.661'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6176'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.370'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.687'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.312'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.097.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8664'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.63%  '
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.24'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.082'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05228'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4282'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.70%  '
$ws.Range("E51").Value = '  +2.72%  '
